$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Product ID 159 -> 13928, Packing Unit ID 3 -> 2, Cart Rules 18 -> 2
$ws.Range("A2").Value = 13928
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2

# Row 3: Product ID 159 -> 13928, Packing Unit ID stays 1, Cart Rules 3 -> 2
$ws.Range("A3").Value = 13928
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2

# Row 4: Product ID 159 -> 5151, Packing Unit ID 16 -> 2, Cart Rules 4 -> 2
$ws.Range("A4").Value = 5151
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
